# error fixed / beta_test.py / v.11.0
#
# Ticket log sheet gets a blank separator row under the header, the
# trailing two empty "Hora de Reparacion" / "Tiempo de Reparacion"
# placeholder cells on the last existing row are cleared out, and four
# new incident rows from 2024-05-17 are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a genuinely blank row 2 (a spacer under the header row) ---
# Touching OutlineLevel forces the row to be materialised in the sheet
# without stamping any cell/style data into it, so it stays empty.
$ws.Rows(2).OutlineLevel = 0

# --- 2. Row 55 no longer carries the two empty trailing cells ---
$ws.Range("H55:I55").ClearContents()

# --- 3. Append the new 2024-05-17 incident rows (56-59) ---
# Column A holds a literal "yyyy-mm-dd" looking string, not a real date,
# so it is entered with a leading apostrophe to force text and keep the
# value exactly as typed instead of Excel re-interpreting it as a date.
# The apostrophe leaves a "quote prefix" style behind, so the cell style
# is put back to Normal right after, once the text value is locked in.

# Row 56
$ws.Range("A56").Value = "'2024-05-17"
$ws.Range("A56").Style = "Normal"
$ws.Range("B56").Value = "09:44:52"
$ws.Range("C56").Value = "-"
$ws.Range("D56").Value = "-"
$ws.Range("E56").Value = "Tornillo atascado"
$ws.Range("F56").Value = "-"
$ws.Range("G56").Value = "-"

# Row 57
$ws.Range("A57").Value = "'2024-05-17"
$ws.Range("A57").Style = "Normal"
$ws.Range("B57").Value = "09:52:45"
$ws.Range("C57").Value = "-"
$ws.Range("D57").Value = "-"
$ws.Range("E57").Value = "Etiquetadora"
$ws.Range("F57").Value = "-"
$ws.Range("G57").Value = "-"
$ws.Range("H57").Value = "09:52:51"
$ws.Range("I57").Value = "0:00:06"

# Row 58
$ws.Range("A58").Value = "'2024-05-17"
$ws.Range("A58").Style = "Normal"
$ws.Range("B58").Value = "09:53:01"
$ws.Range("C58").Value = "-"
$ws.Range("D58").Value = "-"
$ws.Range("E58").Value = "Screw K30 no lo detecta puesto"
$ws.Range("F58").Value = "-"
$ws.Range("G58").Value = "-"
$ws.Range("H58").Value = "09:53:09"
$ws.Range("I58").Value = "0:00:08"

# Row 59
$ws.Range("A59").Value = "'2024-05-17"
$ws.Range("A59").Style = "Normal"
$ws.Range("B59").Value = "09:53:05"
$ws.Range("C59").Value = "-"
$ws.Range("D59").Value = "-"
$ws.Range("E59").Value = "No lee QR"
$ws.Range("F59").Value = "-"
$ws.Range("G59").Value = "-"
$ws.Range("H59").Value = "10:24:34"
$ws.Range("I59").Value = "0:31:29"

Write-Output ("UsedRange after edit: " + $ws.UsedRange.Address())
